$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 251; existing rows 251-348 shift down to 252-349.
$ws.Rows(251).Insert()

# Populate the new row 251 with the new price-report record.
$ws.Cells.Item(251, 1).Value = 4
$ws.Cells.Item(251, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(251, 3).Value = "Los Lagos"
$ws.Cells.Item(251, 4).Value = 44795
$ws.Cells.Item(251, 5).Value = 10
$ws.Cells.Item(251, 6).Value = 100112045
$ws.Cells.Item(251, 7).Value = "Zapallo"
$ws.Cells.Item(251, 8).Value = "Paine"
$ws.Cells.Item(251, 9).Value = "1a (guarda)"
$ws.Cells.Item(251, 10).Value = 700
$ws.Cells.Item(251, 11).Value = 550
$ws.Cells.Item(251, 12).Value = 650
$ws.Cells.Item(251, 13).Value = 600
$ws.Cells.Item(251, 14).Value = "`$/kilo (volumen en unidades)"
$ws.Cells.Item(251, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(251, 16).Value = 600
$ws.Cells.Item(251, 17).Value = 1
$ws.Cells.Item(251, 18).Value = "Hortaliza"
